$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the sample id shared string E7760 -> E7420 across the whole
#    G2:G27 block in one shot so the shared-string slot is edited in place
#    instead of forking a brand new string for a single cell.
$ws.Range("G2:G27").Value = "E7420"

# 2. Give that same G2:G27 block its own font (Arial 11, black) which forces
#    a new font + cell style to be recorded, matching the "further cleaning"
#    metadata split between the G and H columns.
$ws.Range("G2:G27").Font.Name = "Arial"
$ws.Range("G2:G27").Font.Size = 11
$ws.Range("G2:G27").Font.Color = 0

# 3. Turn the literal FALSE booleans in H2:H27 into real =FALSE() formulas
#    (still displaying via the TRUE/FALSE custom number format already on
#    those cells) instead of hard-coded boolean literals.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# 4. Move the live selection from H2:H27 to G2:G27 (active cell G2).
$ws.Range("G2:G27").Select() | Out-Null
